# Update the embedded build-timestamp in the "version" strings that were
# stamped into the workbook at generation time.
#
# Old build timestamp: February 03 2026 17.29.55 EST
# New build timestamp: February 03 2026 18.05.36 EST

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

# --- "About" sheet -------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$a2text = [string]$about.Range("A2").Value()
$about.Range("A2").Value = $a2text.Replace($oldStamp, $newStamp)

$a6text = [string]$about.Range("A6").Value()
$about.Range("A6").Value = $a6text.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet ------------------------------
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 10; $row++) {
    $cell = $data.Range("S$row")
    $stext = [string]$cell.Value()
    $cell.Value = $stext.Replace($oldStamp, $newStamp)
}
